$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing date cell (A19) down to the new row's date cell (A20)
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row's values
$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 0.8976398032236155
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 0.7456737245741252
